$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: clear E7/F7 (were text "NA"), set G7:AC7 to numeric values
$ws.Range("E7").ClearContents()
$ws.Range("F7").ClearContents()
$ws.Range("G7").Value = 4.8568
$ws.Range("H7").Value = 4.6556
$ws.Range("I7").Value = 4.9722
$ws.Range("J7").Value = 5.0324
$ws.Range("K7").Value = 5.9562
$ws.Range("L7").Value = 6.2697
$ws.Range("M7").Value = 5.792000000000001
$ws.Range("N7").Value = 5.8331
$ws.Range("O7").Value = 6.0923
$ws.Range("P7").Value = 6.1935
$ws.Range("Q7").Value = 6.6242
$ws.Range("R7").Value = 6.4175
$ws.Range("S7").Value = 7.0301
$ws.Range("T7").Value = 7.4167
$ws.Range("U7").Value = 8.1306
$ws.Range("V7").Value = 8.394
$ws.Range("W7").Value = 9.2747
$ws.Range("X7").Value = 9.6435
$ws.Range("Y7").Value = 9.2503
$ws.Range("Z7").Value = 9.3676
$ws.Range("AA7").Value = 9.4319
$ws.Range("AB7").Value = 9.3503
$ws.Range("AC7").Value = 9.7046

# Row 19: set E19:AC19 to numeric values (replacing text placeholders)
$ws.Range("E19").Value = 187.41757311
$ws.Range("F19").Value = 195.34687434
$ws.Range("G19").Value = 203.2442
$ws.Range("H19").Value = 209.0379
$ws.Range("I19").Value = 217.4466
$ws.Range("J19").Value = 225.6018
$ws.Range("K19").Value = 233.8428
$ws.Range("L19").Value = 248.0793
$ws.Range("M19").Value = 256.4244
$ws.Range("N19").Value = 259.9155
$ws.Range("O19").Value = 266.7437
$ws.Range("P19").Value = 281.2512
$ws.Range("Q19").Value = 295.5392
$ws.Range("R19").Value = 309.4553
$ws.Range("S19").Value = 327.372
$ws.Range("T19").Value = 346.1355
$ws.Range("U19").Value = 344.1078
$ws.Range("V19").Value = 355.3527
$ws.Range("W19").Value = 371.8061
$ws.Range("X19").Value = 382.8646
$ws.Range("Y19").Value = 387.4156
$ws.Range("Z19").Value = 396.3587
$ws.Range("AA19").Value = 402.9666
$ws.Range("AB19").Value = 414.9285
$ws.Range("AC19").Value = 428.2318

# Row 21: set E21:AC21 to numeric values (replacing text placeholders)
$ws.Range("E21").Value = 187.88957311
$ws.Range("F21").Value = 196.93887434
$ws.Range("G21").Value = 204.2855
$ws.Range("H21").Value = 209.841
$ws.Range("I21").Value = 218.2251
$ws.Range("J21").Value = 226.7279
$ws.Range("K21").Value = 234.9449
$ws.Range("L21").Value = 251.415
$ws.Range("M21").Value = 256.8944
$ws.Range("N21").Value = 259.9237
$ws.Range("O21").Value = 267.9341
$ws.Range("P21").Value = 284.6651000000001
$ws.Range("Q21").Value = 300.0612
$ws.Range("R21").Value = 314.4829
$ws.Range("S21").Value = 331.5511999999999
$ws.Range("T21").Value = 352.3525
$ws.Range("U21").Value = 340.8889
$ws.Range("V21").Value = 358.5406
$ws.Range("W21").Value = 377.1446999999999
$ws.Range("X21").Value = 385.1104
$ws.Range("Y21").Value = 387.4671
$ws.Range("Z21").Value = 397.5216
$ws.Range("AA21").Value = 404.713
$ws.Range("AB21").Value = 417.0524
$ws.Range("AC21").Value = 434.1267

# Row 23: set E23:AC23 to numeric values (replacing text placeholders)
$ws.Range("E23").Value = 296.79758485
$ws.Range("F23").Value = 316.80301925
$ws.Range("G23").Value = 330.1708
$ws.Range("H23").Value = 340.3626
$ws.Range("I23").Value = 362.6118
$ws.Range("J23").Value = 375.3879
$ws.Range("K23").Value = 391.0022
$ws.Range("L23").Value = 436.9671
$ws.Range("M23").Value = 445.6583
$ws.Range("N23").Value = 453.3351
$ws.Range("O23").Value = 461.8117
$ws.Range("P23").Value = 494.8678000000001
$ws.Range("Q23").Value = 529.0657
$ws.Range("R23").Value = 561.8024999999999
$ws.Range("S23").Value = 598.6718
$ws.Range("T23").Value = 634.5288
$ws.Range("U23").Value = 582.6283
$ws.Range("V23").Value = 637.6542999999999
$ws.Range("W23").Value = 686.6302
$ws.Range("X23").Value = 704.0455999999999
$ws.Range("Y23").Value = 708.1276
$ws.Range("Z23").Value = 728.2580999999999
$ws.Range("AA23").Value = 736.2046
$ws.Range("AB23").Value = 767.6118
$ws.Range("AC23").Value = 805.9832000000001

# Row 25: set E25:AC25 to numeric values (replacing text placeholders)
$ws.Range("E25").Value = 41.97260713
$ws.Range("F25").Value = 44.26296716
$ws.Range("G25").Value = 46.1844
$ws.Range("H25").Value = 46.7961
$ws.Range("I25").Value = 50.3804
$ws.Range("J25").Value = 52.5643
$ws.Range("K25").Value = 55.4703
$ws.Range("L25").Value = 61.4099
$ws.Range("M25").Value = 60.0472
$ws.Range("N25").Value = 56.8411
$ws.Range("O25").Value = 58.9599
$ws.Range("P25").Value = 67.1944
$ws.Range("Q25").Value = 73.5383
$ws.Range("R25").Value = 78.0205
$ws.Range("S25").Value = 84.2734
$ws.Range("T25").Value = 91.1028
$ws.Range("U25").Value = 75.56139999999999
$ws.Range("V25").Value = 82.8488
$ws.Range("W25").Value = 90.948
$ws.Range("X25").Value = 89.91760000000001
$ws.Range("Y25").Value = 86.99520000000001
$ws.Range("Z25").Value = 93.2749
$ws.Range("AA25").Value = 96.73729999999999
$ws.Range("AB25").Value = 101.238
$ws.Range("AC25").Value = 107.6641

# Minor floating-point recalculation fix
$ws.Range("E12").Value = 2.734769149999998

